# Weekly update: insert a new data row at the top of the data table (row 4)
# pushing all existing historical rows down by one, and populate the new
# row with the latest week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 4 (shifts rows 4..64 to 5..65)
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest weekly record
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112026
$ws.Range("G4").Value = "Haba"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 7250
$ws.Range("N4").Value = "$/saco 25 kilos"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 290
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
